$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths
# Note: the stored OOXML column "width" attribute = COM ColumnWidth + 5/6 (0.8333...)
# so we subtract that offset here to land on the exact target stored widths.
$ws.Columns.Item(2).ColumnWidth = 121 - 5/6
$ws.Columns.Item(3).ColumnWidth = 25 - 5/6
$ws.Columns.Item(5).ColumnWidth = 11 - 5/6
$ws.Columns.Item(6).ColumnWidth = 11 - 5/6
$ws.Columns.Item(8).ColumnWidth = 11 - 5/6
$ws.Columns.Item(9).ColumnWidth = 11 - 5/6
$ws.Columns.Item(10).ColumnWidth = 22 - 5/6

# Update row 2 cell values
$ws.Range("B2").Value = "Updating Configuration Bills Receivable Config Updating Column File Directory , File Name , Append To File (Y/N) , None"
$ws.Range("C2").Value = "Bills Receivable Config"
$ws.Range("D2").Value = "File Directory"
$ws.Range("E2").Value = "test1"
$ws.Range("F2").Value = "test1"
$ws.Range("G2").Value = "File Name"
$ws.Range("H2").Value = "test2"
$ws.Range("I2").Value = "test2"
$ws.Range("J2").Value = "Append To File (Y/N)"
$ws.Range("K2").Value = "test3"
$ws.Range("L2").Value = "test3"
